$d = $word.ActiveDocument

$replacements = @(
    @{old = "291×9=2619"; new = "542×8=4336"},
    @{old = "247×8=1976"; new = "785×5=3925"},
    @{old = "714×2=1428"; new = "661×9=5949"},
    @{old = "977×3=2931"; new = "665×6=3990"},
    @{old = "218×6=1308"; new = "102×7=714"},
    @{old = "612×5=3060"; new = "782×4=3128"},
    @{old = "683×2=1366"; new = "518×5=2590"},
    @{old = "560×4=2240"; new = "353×5=1765"},
    @{old = "767×9=6903"; new = "595×7=4165"},
    @{old = "838×6=5028"; new = "588×9=5292"},
    @{old = "865×8=6920"; new = "245×5=1225"},
    @{old = "113×8=904";  new = "436×8=3488"},
    @{old = "424×7=2968"; new = "112×5=560"},
    @{old = "670×5=3350"; new = "866×9=7794"},
    @{old = "244×2=488";  new = "457×4=1828"},
    @{old = "715×8=5720"; new = "481×7=3367"},
    @{old = "909×8=7272"; new = "729×2=1458"},
    @{old = "571×3=1713"; new = "589×9=5301"},
    @{old = "746×6=4476"; new = "689×4=2756"},
    @{old = "546×9=4914"; new = "772×6=4632"},
    @{old = "169×3=507";  new = "262×7=1834"},
    @{old = "153×3=459";  new = "312×9=2808"},
    @{old = "660×2=1320"; new = "770×6=4620"},
    @{old = "576×4=2304"; new = "190×5=950"},
    @{old = "478×3=1434"; new = "666×6=3996"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
